$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows.Item(13).Insert()

# New row 13: Docentes responsaveis data (no label, B/C only)
$ws.Range("B13").Value = '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Range("C13").Value = '5840942 - Marco Aurélio Kondracki de Alcântara'

# Row 10 (Objetivos:): replace placeholder text with real PT objectives text
$ws.Range("B10").Value = 'Fornecer aos discentes conhecimentos e informações básicas sobre técnicas e práticas de recuperação de áreas degradadas. Essas informações dizem respeito a legislação e normas vigentes, investigação, caracterização geotécnica, implementação de planos de recuperação e monitoramento.'
$ws.Range("C10").Value = 'Fornecer aos discentes conhecimentos e informações básicas sobre técnicas e práticas de recuperação de áreas degradadas. Essas informações dizem respeito a legislação e normas vigentes, investigação, caracterização geotécnica, implementação de planos de recuperação e monitoramento.'

# Row 14 (Programa resumido:): replace "Semestral" placeholder with real summary text
$ws.Range("B14").Value = 'Introdução e conceitos; identificação do problema: tipos de áreas; legislação e normas; geoindicadores de degradação; técnicas de recuperação de áreas degradadas; implementação de planos de recuperação; monitoramento.'
$ws.Range("C14").Value = 'Introdução e conceitos; identificação do problema: tipos de áreas; legislação e normas; geoindicadores de degradação; técnicas de recuperação de áreas degradadas; implementação de planos de recuperação; monitoramento.'

# Row 16 (Programa:): replace date placeholder with real program text
$ws.Range("B16").Value = 'Degradação e recuperação ambiental; geoindicadores de degradação; legislação e normas aplicadas à recuperação de áreas degradadas; aspectos e níveis de recuperação; tipos de áreas degradadas; técnicas e medidas de recuperação de áreas degradadas; critérios para a seleção de alternativas; implementação de planos de recuperação; monitoramento; exemplos de recuperação de áreas degradadas.'
$ws.Range("C16").Value = 'Degradação e recuperação ambiental; geoindicadores de degradação; legislação e normas aplicadas à recuperação de áreas degradadas; aspectos e níveis de recuperação; tipos de áreas degradadas; técnicas e medidas de recuperação de áreas degradadas; critérios para a seleção de alternativas; implementação de planos de recuperação; monitoramento; exemplos de recuperação de áreas degradadas.'

# Row 19 (Metodo:): replace placeholder with real method text
$ws.Range("B19").Value = 'Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("C19").Value = 'Aulas teóricas e práticas, trabalhos de campo e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios.'

# Row 20 (Criterio:): replace with real criteria text
$ws.Range("B20").Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("C20").Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'

# Row 21 (Norma de recuperacao:): replace with real text
$ws.Range("B21").Value = 'Provas e/ou exercícios dirigidos.'
$ws.Range("C21").Value = 'Provas e/ou exercícios dirigidos.'

# Row 22 (Bibliografia:): replace with full bibliography text
$ws.Range("B22").Value = 'Bibliografia básica:Barrow, C.J. Land Degradation Cambridge University Press, 1991.Berger, A.R. The geoindicator concept and its application: An introduction. In: Berger, A.R. & Iams, W.J. (EDTS) Assessing Rapid Environmental Geoindicators: Changers in Earth Systems. pp.: 1-14 Balkema, Rotterdam, 1996.Blaikie, P. & Brookfield, H. Land degradation and society. London Methuen, 1987Brunsden, D. and Moore, R. Engineering geomorphology on the coast: lessons from West Dorset. Geomorphology 31: 391-409, 1999.CALIJURI, M.C.; CUNHA, D.G.F. Engenharia Ambiental. Conceitos, Tecnologia e Gestão. Rio de Janeiro, Elsevier, 2013.Dahlberg, A.C. Interpretations of environmetal change and diversity: A critical approach to indications of degradation - The case of Kalakamate, Northeast Botswana. Land degradation & Development, 11: 549-562, 2000.DIAS, L.E; GRIFFTITH,J.J. Conceituação e Caracterização de Áreas Degradadas. In: DIAS, L.E; MELLO, J.W.V (orgs). Recuperação de Áreas Degradadas. Universidade Federal de Viçosa. Sociedade Brasileira de Recuperação de Áreas Degradadas, 1998.Duque, M.J.F., Pedroza, J., Ciez, A., Sanz, M.A. & Carrasco, R.M.. A geomorphical design for the rehabilitation of an abandoned sand quarry in central Spain. Landscape and urban planning, 42: 1-14, 1998.GUERRA, A. J. T.; ARAUJO, G., ALMEIDA, J. R. Gestão Ambiental De Áreas Degradadas. Rio de Janeiro : Bertrand Brasil, 2007.Marchetti, M. & Panizza, M. Geomorphology and Environmental Impact Assesssment: A case study in Moema (Dolomites - Italy). In: Marchetti, M & Pinas, V. (EDS). Geomorphology and Environemental Impact Assessements pp: 71-82, Balkema, 2001.MARTINS, S.V. Recuperação de Áreas Degradadas. Universidade Federal de Viçosa. Viçosa, 2013.Neimanis, U. & kerr, A. Developing national environmental indicators. In: Berger, A.R. & Iams, W.J. (EDTS) Assessing rapid environmental geoindicators: changes in earth systems. 1996.SANCHEZ, L.E. Desengenharia: o passive ambiental na desativaçao de empreendimentos industriais. São Paulo, EDUSP, 2001.SANCHEZ, L.E. Avaliação de Impacto Ambiental. São Paulo, Oficina de Textos, 2006.Bibliografia complementar:Berger, A.R. Assessing Rapid Environmetal Change Using Geoindicators. Environmetal Geology, 32, n. 1, 36-44, 1997.Fao. A provisional methodology for soil degradation assessment. FAO. Rome, 1979,Lindskog, P. and Tengberg, A. Land degradation, Natural resources and local knowledge in the Sahel zone of Burkina Faso. Geojournal, 33, 365-375, 1994.Morton, R. A. Geoindicators of coastal wet land and shorelines. In: berger, A.R. & Iams, W.J. (EDTS) Assessing Rapid Environmental Geoindicators: Changes inEarth Systems. pp: 207-232, 1996.Murthy, R.C. Rao, Y. R. and Inamdar, A.B. Integrated coastal management of Mumbai Metropolitan Region. Ocean & Coastal Management 44: 355-369, 2001.'
$ws.Range("C22").Value = 'Bibliografia básica:Barrow, C.J. Land Degradation Cambridge University Press, 1991.Berger, A.R. The geoindicator concept and its application: An introduction. In: Berger, A.R. & Iams, W.J. (EDTS) Assessing Rapid Environmental Geoindicators: Changers in Earth Systems. pp.: 1-14 Balkema, Rotterdam, 1996.Blaikie, P. & Brookfield, H. Land degradation and society. London Methuen, 1987Brunsden, D. and Moore, R. Engineering geomorphology on the coast: lessons from West Dorset. Geomorphology 31: 391-409, 1999.CALIJURI, M.C.; CUNHA, D.G.F. Engenharia Ambiental. Conceitos, Tecnologia e Gestão. Rio de Janeiro, Elsevier, 2013.Dahlberg, A.C. Interpretations of environmetal change and diversity: A critical approach to indications of degradation - The case of Kalakamate, Northeast Botswana. Land degradation & Development, 11: 549-562, 2000.DIAS, L.E; GRIFFTITH,J.J. Conceituação e Caracterização de Áreas Degradadas. In: DIAS, L.E; MELLO, J.W.V (orgs). Recuperação de Áreas Degradadas. Universidade Federal de Viçosa. Sociedade Brasileira de Recuperação de Áreas Degradadas, 1998.Duque, M.J.F., Pedroza, J., Ciez, A., Sanz, M.A. & Carrasco, R.M.. A geomorphical design for the rehabilitation of an abandoned sand quarry in central Spain. Landscape and urban planning, 42: 1-14, 1998.GUERRA, A. J. T.; ARAUJO, G., ALMEIDA, J. R. Gestão Ambiental De Áreas Degradadas. Rio de Janeiro : Bertrand Brasil, 2007.Marchetti, M. & Panizza, M. Geomorphology and Environmental Impact Assesssment: A case study in Moema (Dolomites - Italy). In: Marchetti, M & Pinas, V. (EDS). Geomorphology and Environemental Impact Assessements pp: 71-82, Balkema, 2001.MARTINS, S.V. Recuperação de Áreas Degradadas. Universidade Federal de Viçosa. Viçosa, 2013.Neimanis, U. & kerr, A. Developing national environmental indicators. In: Berger, A.R. & Iams, W.J. (EDTS) Assessing rapid environmental geoindicators: changes in earth systems. 1996.SANCHEZ, L.E. Desengenharia: o passive ambiental na desativaçao de empreendimentos industriais. São Paulo, EDUSP, 2001.SANCHEZ, L.E. Avaliação de Impacto Ambiental. São Paulo, Oficina de Textos, 2006.Bibliografia complementar:Berger, A.R. Assessing Rapid Environmetal Change Using Geoindicators. Environmetal Geology, 32, n. 1, 36-44, 1997.Fao. A provisional methodology for soil degradation assessment. FAO. Rome, 1979,Lindskog, P. and Tengberg, A. Land degradation, Natural resources and local knowledge in the Sahel zone of Burkina Faso. Geojournal, 33, 365-375, 1994.Morton, R. A. Geoindicators of coastal wet land and shorelines. In: berger, A.R. & Iams, W.J. (EDTS) Assessing Rapid Environmental Geoindicators: Changes inEarth Systems. pp: 207-232, 1996.Murthy, R.C. Rao, Y. R. and Inamdar, A.B. Integrated coastal management of Mumbai Metropolitan Region. Ocean & Coastal Management 44: 355-369, 2001.'

